$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = 1.3
$ws.Range("F12").Value = 1.45
$ws.Range("F13").Value = 1.6
$ws.Range("F14").Value = 1.6
$ws.Range("F15").Value = 1.6

$ws.Range("F12").Select()
